$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 3 (MAT143 Brief Calculus) - picture editId change only
$tbl.Rows.Item(3).Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML('<w:p w14:paraId="02B78FEC" w14:textId="37C126B1" w:rsidR="00C11830" w:rsidRDefault="00C11830" w:rsidP="005C117A"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="36273A8B" wp14:editId="0934C163"><wp:extent cx="438411" cy="408241"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1093558619" name="Picture 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 3"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="459821" cy="428178"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')

# Row 5 (STA311 Statistical Comp. and Data Mgmt.) - picture editId change only
$tbl.Rows.Item(5).Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML('<w:p w14:paraId="5B2422E9" w14:textId="3D162B55" w:rsidR="00C11830" w:rsidRDefault="00C11830" w:rsidP="005C117A"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1689430F" wp14:editId="7AD2EF56"><wp:extent cx="463550" cy="406562"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1706856057" name="Picture 1706856057" descr="A logo with text and globe&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1329169576" name="Picture 1" descr="A logo with text and globe&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="490040" cy="429795"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')

# Row 8 (STA501 Methods of Applied Statistics) - picture editId change only
$tbl.Rows.Item(8).Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML('<w:p w14:paraId="2EBD826D" w14:textId="79001262" w:rsidR="00C11830" w:rsidRDefault="00C11830" w:rsidP="005C117A"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="362C5140" wp14:editId="1569118C"><wp:extent cx="455348" cy="407096"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1161971414" name="Picture 5"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 9"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId16" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="477358" cy="426774"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')

# Row 12 (MAT125 Probability and Statistics) - remove centering, reposition + rewrap picture
$tbl.Rows.Item(12).Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML('<w:p w14:paraId="24D00F9F" w14:textId="48E0DFBF" w:rsidR="00E03585" w:rsidRDefault="00E5610C" w:rsidP="005C117A"><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251658240" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="31AEF146" wp14:editId="3374B858"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>80645</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>63500</wp:posOffset></wp:positionV><wp:extent cx="325924" cy="286487"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapThrough wrapText="bothSides"><wp:wrapPolygon edited="0"><wp:start x="0" y="2874"/><wp:lineTo x="0" y="18678"/><wp:lineTo x="20211" y="18678"/><wp:lineTo x="20211" y="2874"/><wp:lineTo x="0" y="2874"/></wp:wrapPolygon></wp:wrapThrough><wp:docPr id="1011079144" name="Picture 1" descr="A black and white text with a green stripe&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1011079144" name="Picture 1" descr="A black and white text with a green stripe&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId24" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="325924" cy="286487"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>')
$tbl.Rows.Item(12).Height = 30.1

# Row 13 (STA506 Mathematical Statistics II) - reposition + rewrap picture
$tbl.Rows.Item(13).Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML('<w:p w14:paraId="187082AE" w14:textId="6B1F5185" w:rsidR="00E03585" w:rsidRDefault="004673EE" w:rsidP="005C117A"><w:pPr><w:jc w:val="center"/><w:rPr><w:noProof/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251660288" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7A70C1EF" wp14:editId="2D259244"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>93345</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>12700</wp:posOffset></wp:positionV><wp:extent cx="325924" cy="286487"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapThrough wrapText="bothSides"><wp:wrapPolygon edited="0"><wp:start x="0" y="2874"/><wp:lineTo x="0" y="18678"/><wp:lineTo x="20211" y="18678"/><wp:lineTo x="20211" y="2874"/><wp:lineTo x="0" y="2874"/></wp:wrapPolygon></wp:wrapThrough><wp:docPr id="1058419839" name="Picture 1058419839" descr="A black and white text with a green stripe&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1011079144" name="Picture 1" descr="A black and white text with a green stripe&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId24" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="325924" cy="286487"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>')
$tbl.Rows.Item(13).Height = 29.65

Write-Output "all edits applied"
